$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to remain plain text so values like "1.00" or
# "58.149.65" are not reinterpreted as numbers/dates by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = "58.203.39"
$ws.Range("E2").Value = "  -0.23%  "

# Row 3
$ws.Range("D3").Value = "2.521.67"
$ws.Range("E3").Value = "  +1.74%  "

# Row 4
$ws.Range("D4").Value = "0.998"
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").Value = "520.93"
$ws.Range("E5").Value = "  +0.02%  "

# Row 6
$ws.Range("D6").Value = "133.02"
$ws.Range("E6").Value = "  -1.29%  "

# Row 7
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  +0.01%  "

# Row 8
$ws.Range("D8").Value = "0.561"
$ws.Range("E8").Value = "  +0.34%  "

# Row 9
$ws.Range("D9").Value = "2.519.43"

# Row 10
$ws.Range("D10").Value = "0.0975"
$ws.Range("E10").Value = "  -0.94%  "

# Row 11
$ws.Range("E11").Value = "  -1.52%  "

# Row 12
$ws.Range("E12").Value = "  -3.40%  "

# Row 13
$ws.Range("E13").Value = "  -2.66%  "

# Row 14
$ws.Range("D14").Value = "2.964.24"
$ws.Range("E14").Value = "  +1.59%  "

# Row 15
$ws.Range("D15").Value = "58.205.60"

# Row 16
$ws.Range("D16").Value = "22.06"
$ws.Range("E16").Value = "  -0.52%  "

# Row 17
$ws.Range("D17").Value = "0.0000135"
$ws.Range("E17").Value = "  -0.59%  "

# Row 18
$ws.Range("D18").Value = "2.507.25"
$ws.Range("E18").Value = "  +0.98%  "

# Row 19
$ws.Range("D19").Value = "10.64"
$ws.Range("E19").Value = "  -0.51%  "

# Row 20
$ws.Range("D20").Value = "321.74"
$ws.Range("E20").Value = "  +0.13%  "

# Row 21
$ws.Range("E21").Value = "  -1.07%  "

# Row 22
$ws.Range("E22").Value = "  +7.16%  "

# Row 23
$ws.Range("D23").Value = "1.00"
$ws.Range("E23").Value = "  +0.06%  "

# Row 24
$ws.Range("D24").Value = "64.61"
$ws.Range("E24").Value = "  +0.16%  "

# Row 25
$ws.Range("E25").Value = "  -1.05%  "

# Row 26
$ws.Range("D26").Value = "0.997"
$ws.Range("E26").Value = "  +0.07%  "

# Row 27
$ws.Range("E27").Value = "  -0.91%  "

# Row 28
$ws.Range("E28").Value = "  -0.42%  "

# Row 29
$ws.Range("E29").Value = "  -0.24%  "

# Row 30
$ws.Range("B30").Value = "Monero"
$ws.Range("C30").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D30").Value = "168.25"
$ws.Range("E30").Value = "  -0.59%  "

# Row 31
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").Value = "1.72"
$ws.Range("E31").Value = "  +1.03%  "

# Row 32
$ws.Range("E32").Value = "  -0.04%  "

# Row 33
$ws.Range("E33").Value = "  -0.22%  "

# Row 34
$ws.Range("E34").Value = "  -0.03%  "

# Row 35
$ws.Range("D35").Value = "0.996"
$ws.Range("E35").Value = "  +0.06%  "

# Row 36
$ws.Range("D36").Value = "18.18"
$ws.Range("E36").Value = "  +0.32%  "

# Row 37
$ws.Range("E37").Value = "  -5.54%  "

# Row 38
$ws.Range("D38").Value = "3.91"
$ws.Range("E38").Value = "  -2.64%  "

# Row 39
$ws.Range("E39").Value = "  +0.64%  "

# Row 40
$ws.Range("D40").Value = "36.38"
$ws.Range("E40").Value = "  -0.67%  "

# Row 41
$ws.Range("D41").Value = "0.767"
$ws.Range("E41").Value = "  -4.13%  "

# Row 42
$ws.Range("D42").Value = "276.32"
$ws.Range("E42").Value = "  +0.39%  "

# Row 43
$ws.Range("D43").Value = "3.44"
$ws.Range("E43").Value = "  -0.64%  "

# Row 44
$ws.Range("D44").Value = "5.02"
$ws.Range("E44").Value = "  -3.18%  "

# Row 45
$ws.Range("D45").Value = "129.42"
$ws.Range("E45").Value = "  +4.17%  "

# Row 46
$ws.Range("D46").Value = "0.598"
$ws.Range("E46").Value = "  -0.04%  "

# Row 47
$ws.Range("D47").Value = "0.0919"
$ws.Range("E47").Value = "  +0.87%  "

# Row 48
$ws.Range("E48").Value = "  +1.79%  "

# Row 49
$ws.Range("D49").Value = "17.68"
$ws.Range("E49").Value = "  -0.55%  "

# Row 50
$ws.Range("E50").Value = "  -0.35%  "

# Row 51
$ws.Range("D51").Value = "16.87"
$ws.Range("E51").Value = "  -1.26%  "
